# Auto-generated edit script applying cryptos.xlsx value updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.517.02"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.376.02"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "310.20"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "104.22"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.79%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.512"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -5.11%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.70%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "35.80"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  +2.15%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0808"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("E13").Value = "  -0.66%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.33%  "
$ws.Range("D15").Value = "2.745.38"
$ws.Range("E15").Value = "  +2.76%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "15.57"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.63%  "
$ws.Range("D17").Value = "2.374.34"
$ws.Range("E17").Value = "  +2.53%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.810"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "43.498.99"
$ws.Range("E19").Value = "  +1.02%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.31"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.92"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").Value = "0.0₃0912"
$ws.Range("E22").Value = "  -1.09%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "68.28"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "240.63"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").Value = "  +0.74%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.60"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.44%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "25.83"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.41%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("E30").Value = "  +9.64%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "36.73"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.29%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "9.48"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.00%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "160.61"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("E34").Value = "  -2.10%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "18.24"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  +5.19%  "
$ws.Range("E38").Value = "  -1.36%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.66"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.60%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0735"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("E41").Value = "  +4.25%  "
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("E43").Value = "  -2.14%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.66"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +15.21%  "
$ws.Range("D45").Value = "2.033.48"
$ws.Range("E45").Value = "  +2.88%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "19.65"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.63%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.12"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "10.57"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +7.51%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "57.82"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.84%  "
$ws.Range("D51").Value = "2.608.29"
$ws.Range("E51").Value = "  +2.74%  "
